# "Added styling to Home page"
# Insert a new header row above the existing data table and style it in
# bold so the sheet reads as name / supply / demand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down one row, keeping all of its values/styles.
$ws.Rows("1:1").Insert() | Out-Null

# Fill in the new header row.
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "supply"
$ws.Range("C1").Value = "demand"

# Bold the header row.
$ws.Range("A1:C1").Font.Bold = $true

# Match the saved selection state (row 1 is frozen as a header, data below it
# is selected).
$ws.Range("A3:C25").Select() | Out-Null
